$wb = $excel.ActiveWorkbook

# Overview sheet: update Latest HO Xliff Generate Date for the dafdef12... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-20 16:56:41"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for the dafdef12... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-20 16:56:37"
$wsZhCn.Range("K3").Value = "2016-08-20 16:56:54"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for the dafdef12... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-20 16:56:41"
$wsDeDe.Range("K3").Value = "2016-08-20 16:57:01"
